$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value2 = '42.796.56'
$ws.Range('E2').Value2 = '  +0.33%  '
$ws.Range('D3').Value2 = '2.310.78'
$ws.Range('E3').Value2 = '  +0.72%  '
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value2 = '  -0.09%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '301.34'
$ws.Range('E5').Value2 = '  -0.25%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '95.18'
$ws.Range('E6').Value2 = '  -0.85%  '
$ws.Range('E7').Value2 = '  +0.08%  '
$ws.Range('E8').Value2 = '  -0.05%  '
$ws.Range('E9').Value2 = '  -0.92%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '34.10'
$ws.Range('E10').Value2 = '  -1.80%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '19.01'
$ws.Range('E11').Value2 = '  +2.49%  '
$ws.Range('E12').Value2 = '  +0.22%  '
$ws.Range('E14').Value2 = '  -1.60%  '
$ws.Range('D15').Value2 = '2.672.47'
$ws.Range('E15').Value2 = '  +0.78%  '
$ws.Range('D16').Value2 = '2.326.26'
$ws.Range('E16').Value2 = '  +0.24%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '0.787'
$ws.Range('E17').Value2 = '  +1.89%  '
$ws.Range('D18').Value2 = '42.701.28'
$ws.Range('E18').Value2 = '  +0.23%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '12.19'
$ws.Range('E19').Value2 = '  -4.77%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '6.13'
$ws.Range('E20').Value2 = '  +2.32%  '
$ws.Range('E21').Value2 = '  -0.18%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '67.68'
$ws.Range('E22').Value2 = '  +0.96%  '
$ws.Range('E23').Value2 = '  +7.29%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '235.02'
$ws.Range('E24').Value2 = '  -0.33%  '
$ws.Range('E25').Value2 = '  -0.02%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '2.42'
$ws.Range('E26').Value2 = '  +1.49%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '24.27'
$ws.Range('E27').Value2 = '  -1.12%  '
$ws.Range('E28').Value2 = '  +15.22%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '166.14'
$ws.Range('E29').Value2 = '  -0.81%  '
$ws.Range('E30').Value2 = '  +1.72%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '32.10'
$ws.Range('E31').Value2 = '  -2.10%  '
$ws.Range('E32').Value2 = '  -0.05%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '5.00'
$ws.Range('E33').Value2 = '  +1.09%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '17.68'
$ws.Range('E34').Value2 = '  -0.40%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '4.44'
$ws.Range('E35').Value2 = '  -0.06%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.0696'
$ws.Range('E36').Value2 = '  +1.59%  '
$ws.Range('E37').Value2 = '  -0.87%  '
$ws.Range('B38').Value2 = 'ARBITRUM'
$ws.Range('C38').Value2 = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '1.77'
$ws.Range('E38').Value2 = '  +2.23%  '
$ws.Range('B39').Value2 = 'Kaspa'
$ws.Range('C39').Value2 = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.100'
$ws.Range('E39').Value2 = '  +0.08%  '
$ws.Range('E40').Value2 = '  +1.15%  '
$ws.Range('E41').Value2 = '  -0.50%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '20.93'
$ws.Range('E42').Value2 = '  +15.42%  '
$ws.Range('D43').Value2 = '1.926.45'
$ws.Range('E43').Value2 = '  -3.28%  '
$ws.Range('E44').Value2 = '  -0.39%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '10.14'
$ws.Range('E45').Value2 = '  +0.03%  '
$ws.Range('E46').Value2 = '  -2.88%  '
$ws.Range('E47').Value2 = '  -0.94%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '2.89'
$ws.Range('E48').Value2 = '  +2.04%  '
$ws.Range('D49').Value2 = '2.540.79'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '53.23'
$ws.Range('E50').Value2 = '  -0.47%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '72.01'
$ws.Range('E51').Value2 = '  +1.81%  '
